$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "220.35", "0.06360") are preserved verbatim as text, matching
# the inlineStr cell type used throughout this sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "26.310.02"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.18"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "220.35"
$ws.Range("E5").Value = "  +0.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5308"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2646"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06360"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "20.93"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07838"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "4.526"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "1.670.23"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "1.895.35"
$ws.Range("E14").Value = "  +0.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5606"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = "65.73"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "26.304.53"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "198.92"
$ws.Range("E21").Value = "  +3.58%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "6.053"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = "146.68"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1214"
$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = "7.236"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value = "16.19"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("E29").Value = "  +2.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05902"
$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value = "1.284"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value = "3.527"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value = "3.316"
$ws.Range("E33").Value = "  +0.92%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D35").Value = "2.827"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9609"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5802"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01616"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "5.961"
$ws.Range("E40").Value = "  +1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075.01"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8573"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.79"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = "1.805.85"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = "58.41"
$ws.Range("E46").Value = "  +2.30%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4413"

$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.078"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₈103"
$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.20%  "
